$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D10").Value = -8.283199999999995
$ws.Range("D12").Value = -7.112599999999998
$ws.Range("D18").Value = -8.796599999999998
$ws.Range("D37").Value = -7.830499999999996
$ws.Range("D55").Value = -8.986399999999998
$ws.Range("D68").Value = -6.865399999999995
$ws.Range("D77").Value = -5.687799999999998
$ws.Range("D78").Value = -7.580200000000004
